$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36, shifting rows 36:93 down to 37:94.
$ws.Rows("36:36").Insert()

# Populate the newly inserted row 36 with the new data point.
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C36").Value = "Ñuble"
$ws.Range("D36").Value = 44915
$ws.Range("E36").Value = 16
$ws.Range("F36").Value = 100112022
$ws.Range("G36").Value = "Arveja Verde"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 50
$ws.Range("K36").Value = 22000
$ws.Range("L36").Value = 22000
$ws.Range("M36").Value = 22000
$ws.Range("N36").Value = "$/saco 25 kilos"
$ws.Range("O36").Value = "Región del Maule"
$ws.Range("P36").Value = 880
$ws.Range("Q36").Value = 25
$ws.Range("R36").Value = "Hortaliza"
